$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows (62 & 63) describing the new "recommender system /
# referral" paragraphs (p24a / p24b), pushing the old row 62 ("Exclusion
# criteria" header) and everything after it down by two rows. ---
$ws.Rows(62).Insert()
$ws.Rows(63).Insert()

# Write the new cells in the same order the cells are laid out left-to-right,
# row-by-row-pairs so the newly created shared strings land at the same
# indices as the authoritative edit (A62, B62, A63, B63, C62, C63).
$ws.Range("A62").Value = "p24a"
$ws.Range("B62").Value = "You may also be compensated for referring the data collection platform to a friend at the rate of 5 EUR per a successful referral. This compensation is only available if you reside in the Netherlands (as is any compensation in this study). A referral is considered successful only if the referred user uses the unique promotion code passed-on by you during the registration, and remains in the study for a period of 1 year. After this period, the reward amount can be requested for payment together with your bank details as explained above. You will be able to review the number of referrals on your dashboard at any time (when logged-in). This referral rewards program will become unavailable for new referrals as soon as the target recruitment is completed and you shall be explicitly informed of this via this website dashboard (when logged-in). This is entirely optional and has no other bearing on the study."

$ws.Range("A63").Value = "p24b"
$ws.Range("B63").Value = "For University of Leiden students performing online tasks, the following reward structure applies. 1 Credit is offered for 30 min of online tasks. This corresponds to performing 2 recommended online tasks for 1 Credit, and a maximum of 4 credits can be obtained by performing 8 recommended tasks. The minimum duration for remaining in the study to receive any credit is 2 months. It is possible to remain in the study for paid or voluntary participation upon termination of the credit period and then general participation commitment applies."

$ws.Range("C62").Value = "U kunt ook een beloning van 5 euro verdienen door iemand succesvol naar de het onderzoek te verwijzen. Deze beloning is alleen beschikbaar als u in Nederland woont (dat geldt voor elke betaling voor deze studie). Een verwijzing geldt alleen als succesvol als de persoon die verwezen is tijdens de registratie de unieke promotiecode invoert die u hebt gegeven, en een jaar lang deel blijft nemen aan het onderzoek. Na deze periode kan de beloning worden opgevraagd onder vermelding van de bankgegevens, zoals hierboven is uitgelegd. U kunt het aantal verwijzingen altijd inzien op uw dashboard (als u ingelogd bent). Zodra het beoogde aantal deelnemers bereikt is zal de mogelijkheid om geld te krijgen voor verwijzingen worden beëindigd. Hierover word u expliciet geïnformeerd via uw dashboard (als u ingelogd bent). Het staat u vrij om wel of niet mee te doen met het verwijzen van deelnemers. Dit heeft verder geen gevolgen voor de studie. "

$ws.Range("C63").Value = "Voor studenten van de Universiteit Leiden die aan online taken deelnemen gelden de volgende beloningsregels. Er wordt 1 SONA credit toegekend voor 30 minuten online-testdeelname. Dit komt overeen met het doorlopen van twee aanbevolen online tests voor 1 credit. Er kunnen maximaal 4 credits worden behaald door 8 aanbevolen taken uit te voeren. U moet minimaal twee maanden aan de studie blijven deelnemen om credits te kunnen verdienen. Na de minimum periode die voor credit beloning geldt is het mogelijk om betaald of vrijwillig door te gaan met de studie. Vanaf dat moment gelden de algemene deelnameregels."

# --- Row-height tidy-up: rows 17 & 48 previously carried stale "tall" heights
# (leftover from longer text); the text removed/shortened ("rt removed before
# ts") now fits in a much shorter row. ---
$ws.Rows(17).RowHeight = 58
$ws.Rows(48).RowHeight = 87

# --- Restore scroll position / selection to where the author was working
# after the edit (near the newly-inserted rows, rather than near the bottom
# of the sheet). ---
[void]$ws.Range("C54").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 49
$win.ScrollColumn = 3
